$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in H1, matching the style/format of the other header cells
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Fill H2:H11 (Iterations=100 block): Controls (rows 2-6) -> 0, MDD (rows 7-11) -> 1
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# Fill H12:H21 (Iterations=200 block): Controls (rows 12-16) -> 0, MDD (rows 17-21) -> 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
